# Auto-generated update of Coeurl Profits workbook
# Applies updated market-price driven values (currentAveragePrice*, LevePrice*, LeveProfit*)
# for rows across multiple job sheets, per scheduled runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1878
$ws.Range("I28").Value = 1128.8462
$ws.Range("K28").Value = 1128.8462
$ws.Range("M28").Value = -643.8462

$ws.Range("H51").Value = 2948.2415
$ws.Range("I51").Value = 2575
$ws.Range("K51").Value = 2575
$ws.Range("M51").Value = -2091

$ws.Range("H86").Value = 8208.166999999999
$ws.Range("I86").Value = 8000
$ws.Range("K86").Value = 8000
$ws.Range("M86").Value = -6877

$ws.Range("H89").Value = 8208.166999999999
$ws.Range("I89").Value = 8000
$ws.Range("K89").Value = 40000
$ws.Range("M89").Value = -34384

$ws.Range("H100").Value = 2430.2632
$ws.Range("I100").Value = 2421.7058
$ws.Range("K100").Value = 2421.7058
$ws.Range("M100").Value = -1880.7058

$ws.Range("H121").Value = 2376.75
$ws.Range("J121").Value = 2376.75
$ws.Range("L121").Value = 7130.25
$ws.Range("N121").Value = -10624.25

$ws.Range("H132").Value = 1499.5156
$ws.Range("I132").Value = 1488.0328
$ws.Range("K132").Value = 4464.0984
$ws.Range("M132").Value = -1934.0984

$ws.Range("H137").Value = 4236.1816
$ws.Range("I137").Value = 5283.1665
$ws.Range("J137").Value = 2979.8
$ws.Range("K137").Value = 15849.4995
$ws.Range("L137").Value = 8939.400000000001
$ws.Range("M137").Value = -13299.4995
$ws.Range("N137").Value = -14039.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1836.7858
$ws.Range("J2").Value = 1971.1111
$ws.Range("L2").Value = 1971.1111
$ws.Range("N2").Value = -2197.1111

$ws.Range("H32").Value = 22796.414
$ws.Range("I32").Value = 21147.713
$ws.Range("J32").Value = 50000
$ws.Range("K32").Value = 21147.713
$ws.Range("L32").Value = 50000
$ws.Range("M32").Value = -20860.713
$ws.Range("N32").Value = -50574

$ws.Range("H61").Value = 254570.38
$ws.Range("I61").Value = 3610.85
$ws.Range("K61").Value = 3610.85
$ws.Range("M61").Value = -3398.85

$ws.Range("H74").Value = 26303.154
$ws.Range("I74").Value = 12894.1
$ws.Range("K74").Value = 12894.1
$ws.Range("M74").Value = -12020.1

$ws.Range("H77").Value = 26303.154
$ws.Range("I77").Value = 12894.1
$ws.Range("K77").Value = 64470.5
$ws.Range("M77").Value = -60102.5

$ws.Range("H102").Value = 1608.9375
$ws.Range("I102").Value = 1608.9375
$ws.Range("K102").Value = 1608.9375
$ws.Range("M102").Value = 13.0625

$ws.Range("H110").Value = 13015.263
$ws.Range("I110").Value = 17097.166
$ws.Range("K110").Value = 17097.166
$ws.Range("M110").Value = -15052.166

$ws.Range("H116").Value = 1836.7858
$ws.Range("J116").Value = 1971.1111
$ws.Range("L116").Value = 1971.1111
$ws.Range("N116").Value = -6559.1111

$ws.Range("H122").Value = 1660.75
$ws.Range("I122").Value = 1693.4
$ws.Range("K122").Value = 5080.200000000001
$ws.Range("M122").Value = -2630.200000000001

$ws.Range("H136").Value = 254570.38
$ws.Range("I136").Value = 3610.85
$ws.Range("K136").Value = 10832.55
$ws.Range("M136").Value = -8282.549999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1836.7858
$ws.Range("J3").Value = 1971.1111
$ws.Range("L3").Value = 1971.1111
$ws.Range("N3").Value = -2199.1111

$ws.Range("H100").Value = 26110.8
$ws.Range("J100").Value = 26110.8
$ws.Range("L100").Value = 26110.8
$ws.Range("N100").Value = -28274.8

$ws.Range("H107").Value = 1067.579
$ws.Range("I107").Value = 1090.8823
$ws.Range("J107").Value = 869.5
$ws.Range("K107").Value = 1090.8823
$ws.Range("L107").Value = 869.5
$ws.Range("M107").Value = 829.1177
$ws.Range("N107").Value = -4709.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3400.3845
$ws.Range("I58").Value = 3400.3845
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 3400.3845
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -3197.3845
$ws.Range("N58").ClearContents()

$ws.Range("H81").Value = 25000
$ws.Range("J81").Value = 25000
$ws.Range("L81").Value = 25000
$ws.Range("N81").Value = -26996

$ws.Range("H84").Value = 25000
$ws.Range("J84").Value = 25000
$ws.Range("L84").Value = 75000
$ws.Range("N84").Value = -84984

$ws.Range("H97").Value = 44333.332
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H133").Value = 104875
$ws.Range("J133").Value = 104875
$ws.Range("L133").Value = 104875
$ws.Range("N133").Value = -109935

$ws.Range("H136").Value = 3400.3845
$ws.Range("I136").Value = 3400.3845
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10201.1535
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7651.1535
$ws.Range("N136").ClearContents()

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 14460.333
$ws.Range("I56").Value = 14460.333
$ws.Range("K56").Value = 14460.333
$ws.Range("M56").Value = -13930.333

$ws.Range("H80").Value = 7900
$ws.Range("J80").Value = 7900
$ws.Range("L80").Value = 23700
$ws.Range("N80").Value = -25572

$ws.Range("H83").Value = 7900
$ws.Range("J83").Value = 7900
$ws.Range("L83").Value = 71100
$ws.Range("N83").Value = -80460

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4314.1665
$ws.Range("I113").Value = 3971.75
$ws.Range("J113").Value = 4999
$ws.Range("K113").Value = 3971.75
$ws.Range("L113").Value = 4999
$ws.Range("M113").Value = -1801.75
$ws.Range("N113").Value = -9339

$ws.Range("H126").Value = 14516.8
$ws.Range("I126").Value = 17789.268
$ws.Range("K126").Value = 53367.804
$ws.Range("M126").Value = -50897.804

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H61").Value = 63599.59
$ws.Range("I61").Value = 87307.75
$ws.Range("K61").Value = 87307.75
$ws.Range("M61").Value = -87105.75

$ws.Range("H99").Value = 8000
$ws.Range("I99").Value = 8000
$ws.Range("K99").Value = 8000
$ws.Range("M99").Value = -5005

$ws.Range("H113").Value = 63599.59
$ws.Range("I113").Value = 87307.75
$ws.Range("K113").Value = 87307.75
$ws.Range("M113").Value = -85137.75

$ws.Range("H132").Value = 4594.967
$ws.Range("I132").Value = 4607.75
$ws.Range("J132").Value = 4580.357
$ws.Range("K132").Value = 13823.25
$ws.Range("L132").Value = 13741.071
$ws.Range("M132").Value = -11293.25
$ws.Range("N132").Value = -18801.071

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4014
$ws.Range("I132").Value = 3757.6316
$ws.Range("K132").Value = 11272.8948
$ws.Range("M132").Value = -8742.8948

$ws.Range("H136").Value = 1882.3823
$ws.Range("I136").Value = 1834.5172
$ws.Range("K136").Value = 5503.5516
$ws.Range("M136").Value = -2953.5516

$ws.Range("H139").Value = 106129.8
$ws.Range("I139").Value = 105325
$ws.Range("J139").Value = 106666.336
$ws.Range("K139").Value = 105325
$ws.Range("L139").Value = 106666.336
$ws.Range("M139").Value = -100185
$ws.Range("N139").Value = -116946.336
